$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The sheet ships protected; temporarily unprotect so the cells below can
# be updated, then re-apply protection afterwards.
$ws.Unprotect()

# Update the confidentiality / as-of-date notice text (cell A11).
$ws.Range("A11").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-06-14 for illustrative purposes only and are subject to change."

# Update the Weight (D) and Percent Change (E) values for rows 2-8.
$ws.Range("D2").Value = 0.5390778285537491
$ws.Range("E2").Value = -0.004606449028640158

$ws.Range("D3").Value = 0.253349841752445
$ws.Range("E3").Value = 0.00885334856490072

$ws.Range("D4").Value = 0.05026613890474265
$ws.Range("E4").Value = -0.004453240969816918

$ws.Range("D5").Value = 0.09880028813735929
$ws.Range("E5").Value = -0.01208348590259967

$ws.Range("D6").Value = 0.02860307835161765
$ws.Range("E6").Value = -0.00876141289311072

$ws.Range("D7").Value = 0.02990282430008643
$ws.Range("E7").Value = -0.002041587901701503

$ws.Range("E8").Value = -0.001969591823621664

# Restore sheet protection.
$ws.Protect()
